# Team Intro.xlsx -- add README / Sheet2 summary table, add Ronn Pang row,
# tweak a couple of totals, and repoint the active sheet/selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$nbsp = [char]0xA0

# ---------------------------------------------------------------------
# Sheet1: new team-member row (Ronn Pang), and updated "Could also do"
# totals (B19 4 -> 5, which ripples the SUM in B21 from 11 -> 12).
# ---------------------------------------------------------------------
$ws1.Range("A15").Value = "Ronn Pang"
$ws1.Range("B15").Value = 5
$ws1.Range("C15").Value = 8
$ws1.Range("D15").Value = 4
$ws1.Range("E15").Value = 4
$ws1.Range("F15").Value = "Programming: Implementation"
$ws1.Range("G15").Value = "I am still on the waitlist"
$ws1.Range("H15").Value = "Programming: Database or Documentation"

$ws1.Range("B19").Value = 5

# ---------------------------------------------------------------------
# Sheet2: a "Name / Wants to do / Total" README-style roll-up, sorted
# alphabetically by Name.
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 25.89
$ws2.Columns.Item(2).ColumnWidth = 39.78
$ws2.Columns.Item(3).ColumnWidth = 25.89

$ws2.Range("A2").Value = "Name"
$ws2.Range("B2").Value = "Wants to do"
$ws2.Range("C2").Value = "Could also do"
$ws2.Range("A2:C2").Style = "40% - Accent1"

$ws2.Range("A3").Value = "Alessandro Power"
$ws2.Range("B3").Value = "Documentation"

$ws2.Range("A4").Value = "Andy Nguyen"
$ws2.Range("B4").Value = "Documentation"

$ws2.Range("A5").Value = "Anhkhoi Vu-Nguyen"
$ws2.Range("B5").Value = "Programming: Backend"
$ws2.Range("C5").Value = "Anything"

$ws2.Range("A6").Value = "Eric Payette"
$ws2.Range("B6").Value = "Programming: Backend"
$ws2.Range("C6").Value = "Documentation, Front End"

$ws2.Range("A7").Value = "Jacqueline Luo"
$ws2.Range("B7").Value = "Documentation"
$ws2.Range("C7").Value = "Programming but limited; Testing"

$ws2.Range("A8").Value = "James Talarico"
$ws2.Range("B8").Value = "Documentation"

$ws2.Range("A9").Value = "Kenny Nguyen"
$ws2.Range("B9").Value = "Documentation"

$ws2.Range("A10").Value = "Laurendy Lam"
$ws2.Range("B10").Value = "Programming: Everything"
$ws2.Range("C10").Value = "Anything"

$ws2.Range("A11").Value = "Michael Mescheder$nbsp"
$ws2.Range("B11").Value = "Programming: Full-Stack"

$ws2.Range("A12").Value = "Piratheeban Annamalai"
$ws2.Range("B12").Value = "Anything"

$ws2.Range("A13").Value = "Pragas Velauthapillai$nbsp"
$ws2.Range("B13").Value = "Documentation$nbsp"
$ws2.Range("C13").Value = "Programming but limited; Testing"

$ws2.Range("A14").Value = "Ronn Pang"
$ws2.Range("B14").Value = "Programming: Implementation"
$ws2.Range("C14").Value = "Programming: Database or Documentation"

# Record a sort (by Name, ascending) over the data rows so the sheet
# carries a sortState the same way Excel would after Data > Sort.
$ws2.Sort.SortFields.Clear()
$ws2.Sort.SortFields.Add($ws2.Range("A3"))
$ws2.Sort.SetRange($ws2.Range("A3:C14"))
$ws2.Sort.Header = 0
$ws2.Sort.Apply()

# ---------------------------------------------------------------------
# View state: Sheet2 becomes the active tab/sheet, scrolled down a bit,
# with the whole table selected; Sheet1's selection moves to H15.
# ---------------------------------------------------------------------
$ws1.Range("H15").Select()

$ws2.Activate()
$ws2.Range("A2:C14").Select()
$excel.ActiveWindow.ScrollRow = 5
